$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 1766.6666
$ws.Range("I4").Value = 1700
$ws.Range("K4").Value = 1700
$ws.Range("M4").Value = -1586
$ws.Range("H43").Value = 595.6667
$ws.Range("I43").Value = 595.6667
$ws.Range("K43").Value = 595.6667
$ws.Range("M43").Value = -526.6667
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("N72").ClearContents()
$ws.Range("H80").Value = 3070.4
$ws.Range("I80").Value = 3784
$ws.Range("J80").Value = 2000
$ws.Range("K80").Value = 11352
$ws.Range("L80").Value = 6000
$ws.Range("M80").Value = -10354
$ws.Range("N80").Value = -7996
$ws.Range("H83").Value = 3070.4
$ws.Range("I83").Value = 3784
$ws.Range("J83").Value = 2000
$ws.Range("K83").Value = 34056
$ws.Range("L83").Value = 18000
$ws.Range("M83").Value = -29064
$ws.Range("N83").Value = -27984
$ws.Range("H98").Value = 1517.0416
$ws.Range("I98").Value = 1295.5294
$ws.Range("J98").Value = 2055
$ws.Range("K98").Value = 1295.5294
$ws.Range("L98").Value = 2055
$ws.Range("M98").Value = 202.4706000000001
$ws.Range("N98").Value = -5051
$ws.Range("H116").Value = 14160
$ws.Range("I116").Value = 51250
$ws.Range("K116").Value = 51250
$ws.Range("M116").Value = -47808
$ws.Range("H122").Value = 1517.0416
$ws.Range("I122").Value = 1295.5294
$ws.Range("J122").Value = 2055
$ws.Range("K122").Value = 3886.5882
$ws.Range("L122").Value = 6165
$ws.Range("M122").Value = -1436.5882
$ws.Range("N122").Value = -11065
$ws.Range("H135").Value = 111112890
$ws.Range("I135").Value = 1984.2
$ws.Range("J135").Value = 250001520
$ws.Range("K135").Value = 17857.8
$ws.Range("L135").Value = 2250013680
$ws.Range("M135").Value = -15322.8
$ws.Range("N135").Value = -2250018750
$ws.Range("H138").Value = 1749.0182
$ws.Range("I138").Value = 1240.0938
$ws.Range("J138").Value = 2457.087
$ws.Range("K138").Value = 3720.2814
$ws.Range("L138").Value = 7371.261
$ws.Range("M138").Value = 1419.7186
$ws.Range("N138").Value = -17651.261
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7351.057
$ws.Range("I32").Value = 5492.857
$ws.Range("K32").Value = 5492.857
$ws.Range("M32").Value = -5205.857
$ws.Range("H45").Value = 1424.6923
$ws.Range("I45").Value = 1018.1429
$ws.Range("K45").Value = 1018.1429
$ws.Range("M45").Value = -641.1429000000001
$ws.Range("H63").Value = 3651
$ws.Range("I63").Value = 3651
$ws.Range("K63").Value = 3651
$ws.Range("M63").Value = -2965
$ws.Range("H66").Value = 3651
$ws.Range("I66").Value = 3651
$ws.Range("K66").Value = 18255
$ws.Range("M66").Value = -14823
$ws.Range("H132").Value = 1746.2727
$ws.Range("I132").Value = 1421
$ws.Range("K132").Value = 4263
$ws.Range("M132").Value = -1733
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1239.375
$ws.Range("I94").Value = 569.1667
$ws.Range("K94").Value = 569.1667
$ws.Range("M94").Value = -118.1667
$ws.Range("H99").Value = 1998
$ws.Range("J99").Value = 1998
$ws.Range("L99").Value = 1998
$ws.Range("N99").Value = -4994
$ws.Range("H107").Value = 699.5909
$ws.Range("I107").Value = 488.85715
$ws.Range("J107").Value = 1068.375
$ws.Range("K107").Value = 488.85715
$ws.Range("L107").Value = 1068.375
$ws.Range("M107").Value = 1431.14285
$ws.Range("N107").Value = -4908.375
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()
$ws.Range("H122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H134").Value = 11607
$ws.Range("I134").Value = 16239
$ws.Range("K134").Value = 48717
$ws.Range("M134").Value = -46182
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 5436046
$ws.Range("I58").Value = 7247061.5
$ws.Range("K58").Value = 7247061.5
$ws.Range("M58").Value = -7246858.5
$ws.Range("H99").Value = 3218.75
$ws.Range("I99").Value = 2150
$ws.Range("K99").Value = 2150
$ws.Range("M99").Value = -652
$ws.Range("H107").Value = 771.05884
$ws.Range("I107").Value = 531
$ws.Range("J107").Value = 1114
$ws.Range("K107").Value = 531
$ws.Range("L107").Value = 1114
$ws.Range("M107").Value = 1389
$ws.Range("N107").Value = -4954
$ws.Range("H126").Value = 3218.75
$ws.Range("I126").Value = 2150
$ws.Range("K126").Value = 6450
$ws.Range("M126").Value = -3980
$ws.Range("H132").Value = 2547.7144
$ws.Range("I132").Value = 1758.875
$ws.Range("J132").Value = 3599.5
$ws.Range("K132").Value = 5276.625
$ws.Range("L132").Value = 10798.5
$ws.Range("M132").Value = -2746.625
$ws.Range("N132").Value = -15858.5
$ws.Range("H134").Value = 2207
$ws.Range("I134").Value = 2140.9167
$ws.Range("K134").Value = 6422.750100000001
$ws.Range("M134").Value = -3887.750100000001
$ws.Range("H136").Value = 5436046
$ws.Range("I136").Value = 7247061.5
$ws.Range("K136").Value = 21741184.5
$ws.Range("M136").Value = -21738634.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 3313.8572
$ws.Range("J22").Value = 3366.1667
$ws.Range("L22").Value = 10098.5001
$ws.Range("N22").Value = -10436.5001
$ws.Range("H27").Value = 3313.8572
$ws.Range("J27").Value = 3366.1667
$ws.Range("L27").Value = 10098.5001
$ws.Range("N27").Value = -10302.5001
$ws.Range("H131").Value = 30455.709
$ws.Range("I131").Value = 630
$ws.Range("J131").Value = 38304.58
$ws.Range("K131").Value = 1890
$ws.Range("L131").Value = 114913.74
$ws.Range("M131").Value = 3150
$ws.Range("N131").Value = -124993.74
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 5424118.5
$ws.Range("J12").Value = 2842002.2
$ws.Range("L12").Value = 2842002.2
$ws.Range("N12").Value = -2842282.2
$ws.Range("H97").Value = 973.8
$ws.Range("I97").Value = 1010.0769
$ws.Range("K97").Value = 1010.0769
$ws.Range("M97").Value = -514.0769
$ws.Range("H126").Value = 4042481.5
$ws.Range("I126").Value = 6175863
$ws.Range("J126").Value = 202395
$ws.Range("K126").Value = 18527589
$ws.Range("L126").Value = 607185
$ws.Range("M126").Value = -18525119
$ws.Range("N126").Value = -612125
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1982.8182
$ws.Range("I46").Value = 1195
$ws.Range("J46").Value = 2433
$ws.Range("K46").Value = 1195
$ws.Range("L46").Value = 2433
$ws.Range("M46").Value = -1007
$ws.Range("N46").Value = -2809
$ws.Range("H68").Value = 1425.5
$ws.Range("I68").Value = 1425.5
$ws.Range("K68").Value = 1425.5
$ws.Range("M68").Value = -676.5
$ws.Range("H71").Value = 1425.5
$ws.Range("I71").Value = 1425.5
$ws.Range("K71").Value = 7127.5
$ws.Range("M71").Value = -3383.5
$ws.Range("H122").Value = 5184.2915
$ws.Range("J122").Value = 6000.625
$ws.Range("L122").Value = 18001.875
$ws.Range("N122").Value = -22901.875
$ws.Range("H132").Value = 2289.4722
$ws.Range("I132").Value = 1459.6154
$ws.Range("K132").Value = 4378.8462
$ws.Range("M132").Value = -1848.8462
$ws.Range("H136").Value = 5853
$ws.Range("I136").Value = 4413.25
$ws.Range("K136").Value = 13239.75
$ws.Range("M136").Value = -10689.75
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()
$ws.Range("H113").Value = 439.0357
$ws.Range("I113").Value = 331
$ws.Range("K113").Value = 993
$ws.Range("M113").Value = 1177
$ws.Range("H122").Value = 36580.816
$ws.Range("I122").Value = 49662.25
$ws.Range("K122").Value = 148986.75
$ws.Range("M122").Value = -146536.75
$ws.Range("H126").Value = 4992.4
$ws.Range("I126").Value = 5904
$ws.Range("K126").Value = 17712
$ws.Range("M126").Value = -15242
$ws.Range("H135").Value = 143357.2
$ws.Range("J135").Value = 143357.2
$ws.Range("L135").Value = 143357.2
$ws.Range("N135").Value = -153497.2
$ws.Range("H136").Value = 21369030
$ws.Range("I136").Value = 32680782
$ws.Range("K136").Value = 98042346
$ws.Range("M136").Value = -98039796
$ws.Range("H141").Value = 85918.57000000001
$ws.Range("J141").Value = 85918.57000000001
$ws.Range("L141").Value = 85918.57000000001
$ws.Range("N141").Value = -96278.57000000001
